# Update election results for row 2 (BRAGANÇA / VIMIOSO)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 58
$ws.Range("I2").Value = 172
$ws.Range("J2").Value = 657
$ws.Range("L2").Value = 183
$ws.Range("M2").Value = 16
$ws.Range("N2").Value = 97
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 71
$ws.Range("T2").Value = 124
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 982
$ws.Range("X2").Value = 959
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 11
